$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BESS capacity value (B2) from 4 to 2 (BESS degradation applied to NPV)
$ws.Range("B2").Value = 2

# Move the active cell selection to B3 (matches author's final cursor position)
$ws.Range("B3").Select()
